$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('X MIPA 2', 'DEFI SETIYOWATI')
    ,@('X MIPA 6', 'LUKE ARINDA FEBRYOLA')
    ,@('X MIPA 6', 'NADIA RAHMA WINDARINI')
    ,@('X MIPA 1', 'EVI BUDI APRIYANI')
    ,@('X MIPA 1', 'LINDA SRI LESTARI')
    ,@('X MIPA 2', 'ALFIRA SELFIANA PUTRI')
    ,@('X MIPA 2', 'POPPY TRI UTAMI')
    ,@('X MIPA 4', 'FAUZAN BILAL MAHARDIKA')
    ,@('X MIPA 4', 'RICO HENDRIAWAN')
    ,@('X MIPA 6', 'LALA ERLINDA')
    ,@('X MIPA 1', 'ANDIKA IRDI  PERMANA')
    ,@('X MIPA 5', 'ALINSYA NICO OVIYANA')
    ,@('X MIPA 6', 'ADELLIA HANDAYANI')
    ,@('X MIPA 3', 'DHIMAS YUDIYATMOKO')
    ,@('X MIPA 3', 'EKA PANJI SATRIA')
    ,@('X MIPA 3', 'GEISKA AYU WULAN ANGGRAENY')
    ,@('X MIPA 3', 'INTAN MUSTIKA SARI')
    ,@('X MIPA 4', 'PUJI LESTARI')
    ,@('X MIPA 6', 'RIZKY ADI NUGROHO')
    ,@('X MIPA 6', 'ANNISA DHEKA CAHYANINGTYAS')
    ,@('X MIPA 5', 'GERA DWIDYA AYU SHELYMAR')
    ,@('X MIPA 3', 'SYIFA AMANDA PUTRI')
    ,@('X MIPA 3', 'APNAMIRA DWI NOOR RIYANTI')
    ,@('X MIPA 3', 'AMELIANA PUTRI')
    ,@('X MIPA 6', 'YOVI SAPUTRA')
    ,@('X MIPA 5', 'FARIDA DWI LESTARI')
    ,@('X MIPA 2', 'INDRA SETIAWAN')
    ,@('X MIPA 2', 'MOHAMAD TOUFIK HIDAYAT')
    ,@('X MIPA 4', 'FAJAR QOMARUDIN')
    ,@('X MIPA 5', 'SEPTA AL AKBAR')
    ,@('X MIPA 5', 'RENDITYA HIKMAL ARYANTARA')
    ,@('X MIPA 5', 'NAJWA ANJELIA RAHMA TANTIA')
    ,@('X MIPA 5', 'FABIZAN ARKIANTO')
    ,@('X MIPA 5', 'DITO OKTA SETYAWAN')
    ,@('X MIPA 6', 'KEVIN NUR SO''IM')
    ,@('X MIPA 5', 'DIA AYU LESTARI')
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
